# The worksheet contains a 3x3 grid of pictures anchored (via oneCellAnchor)
# at rows 0, 20, 40 (0-based row index) and columns 0, 10, 20.
# This edit pushes the two lower rows of pictures further down:
#   row 20 -> row 30
#   row 40 -> row 60
# (the top row of pictures, anchored at row 0, is left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowHeight = $ws.Rows.Item(1).RowHeight

for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)

    # Determine which 0-based anchor row this picture currently sits on,
    # from its current pixel/point Top position.
    $currentAnchorRow = [Math]::Round($shp.Top / $rowHeight)

    if ($currentAnchorRow -eq 20) {
        $shp.Top = 30 * $rowHeight
    }
    elseif ($currentAnchorRow -eq 40) {
        $shp.Top = 60 * $rowHeight
    }
}
